# Insert a new data row before row 52 (shifts rows 52..114 down to 53..115)
# and populate the new row with a "Macroferia Regional de Talca - Apio"
# weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("52:52").Insert()

$ws.Range("A52").Value = 5
$ws.Range("B52").Value = "Macroferia Regional de Talca"
$ws.Range("C52").Value = "Maule"
$ws.Range("D52").Value = 44467
$ws.Range("E52").Value = 7
$ws.Range("F52").Value = 100112017
$ws.Range("G52").Value = "Apio"
$ws.Range("H52").Value = "Americana (o)"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 500
$ws.Range("K52").Value = 9000
$ws.Range("L52").Value = 9000
$ws.Range("M52").Value = 9000
$ws.Range("N52").Value = "`$/docena de matas"
$ws.Range("O52").Value = "Provincia del Elquí"
$ws.Range("P52").Value = 1500
$ws.Range("Q52").Value = 6
$ws.Range("R52").Value = "Hortaliza"
